$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetimes for first data row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-21 12:59:43"
$wsZh.Range("H2").Value = "2016-03-21 13:00:06"

# de-de sheet: update handoff/handback datetimes for first data row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-21 12:59:46"
$wsDe.Range("H2").Value = "2016-03-21 13:00:15"
